$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the B1 header label to match the other "HiRise_ID" header (E1) and make it bold
$ws.Range("B1").Value = "HiRise_ID"
$ws.Range("B1").Font.Bold = $true

# Set the page to print in portrait orientation
$ws.PageSetup.Orientation = 1

# Move the active selection to C4 (as left by the author)
$null = $ws.Range("C4").Select()
